$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# "Invoices to collect" sheet: replace the two open invoices with a new
# pair of invoice numbers, and make the collection amount a plain
# left-aligned number instead of a shared "50" string.
# ---------------------------------------------------------------------
$wsInvoices = $wb.Worksheets.Item("Invoices to collect")
$wsInvoices.Range("A2:B2").NumberFormat = "@"
$wsInvoices.Range("A2").Value = "I-VS004002270"
$wsInvoices.Range("B2").Value = "Full"
$wsInvoices.Range("A3").Value = "I-VS004002271"
$wsInvoices.Range("B3").Value = 25000
$wsInvoices.Range("B3").HorizontalAlignment = -4131

# ---------------------------------------------------------------------
# "Payment methods" sheet: the Cash / CDC / PDC amounts become text
# values and the CDC cheque date moves from 05 June 2018 to 27 June 2018.
# ---------------------------------------------------------------------
$wsPayments = $wb.Worksheets.Item("Payment methods")
$wsPayments.Range("B2").NumberFormat = "@"
$wsPayments.Range("B2").Value = "7842.88"
$wsPayments.Range("B3").NumberFormat = "@"
$wsPayments.Range("B3").Value = "12000"
$wsPayments.Range("F3").Value = "27 June 2018"
$wsPayments.Range("B4").NumberFormat = "@"
$wsPayments.Range("B4").Value = "8000"

# ---------------------------------------------------------------------
# Selections / active sheet: "Invoices to collect" becomes the selected
# tab (with B3 highlighted); "Payment methods" keeps a live selection
# of F4 for when the user returns to it.
# ---------------------------------------------------------------------
$wsPayments.Range("F4").Select()
$wsInvoices.Activate()
$wsInvoices.Range("B3").Select()
